$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.647.47'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.982.54'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.76%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.90'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('E6').Value = '  -4.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.73'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +6.20%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '58.96'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0737'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.953'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.57'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.271.43'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.86'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +10.03%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.978.78'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.12%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '35.556.71'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.92%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '71.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0848'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '233.13'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.58'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +20.91%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.28'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.99%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.06'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.24'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.67%  '
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.89'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.38%  '
$ws.Range('E32').Value = '  -7.69%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0952'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +12.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0596'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.44'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +9.59%  '
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.45'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +8.94%  '
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.87'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0214'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '93.81'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.09'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0909'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.18'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.373.52'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '46.76'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('E51').Value = '  -0.38%  '
